$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Skills line: "Git, Machine Learning" -> "Git, SQL, Machine Learning"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Git, Machine Learning", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Git, SQL, Machine Learning", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) "| UGAHACKS11 : NCR VOYX Track" -> "| UGAHACKS11: NCR VOYX Track"
#    (drop the stray space before the colon)
# ---------------------------------------------------------------------
$rng = $d.Content
$old2 = "| UGAHACKS" + "11 :" + " NCR VOYX Track"
$new2 = "| UGAHACKS11: NCR VOYX Track"
$rng.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Drop "register_order " before "RPC that processes live POS orders..."
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("including an atomic register_order RPC", $true, $false, $false, $false, $false, $true, 1, $false, `
    "including an atomic RPC", 2) | Out-Null

# Clean up the now-orphaned grammar markers around "joins" by re-writing
# the phrase with itself (forces the run / proofErr bookkeeping to be
# rebuilt around the current text).
$rng = $d.Content
$rng.Find.Execute("Bill of Materials joins, and decrementing", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Bill of Materials joins, and decrementing", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Remove the empty "ListParagraph" line right before "Speech Mate ..."
# ---------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "" -and $p.Range.ParagraphStyle.NameLocal -eq "List Paragraph") {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text.StartsWith("Speech Mate")) {
            $p.Range.Delete()
            $found = $true
            break
        }
    }
}
Write-Output ("Removed empty paragraph before Speech Mate: " + $found)

# ---------------------------------------------------------------------
# 5) "| AI-ATL Hackathon (Georgia Tech) | Spring 2025"
#    -> "| AI-ATL Hackathon (Georgia Tech) | " + italic "Fall 2025"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("| AI-ATL Hackathon (Georgia Tech) | Spring 2025", $true, $false, $false, $false, $false, $true, 1, $false, `
    "| AI-ATL Hackathon (Georgia Tech) | Fall 2025", 2) | Out-Null

$rng2 = $d.Content
$rng2.Find.Execute("Fall 2025", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$f = $rng2.Font
$f.Italic = 1
$f.ItalicBi = 1

Write-Output "done"
